# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de), the handback step now
# fills in "Latest Target File" (F) and "Latest Handback File" (G) with the
# same file links already recorded in "Source File Name" (A) / "Latest
# Handoff File" (D), stamps "Latest Handback DateTime" (H) with the time the
# handback finished, and flips the Status text from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is shown (the Overview
# sheet included, since it mirrors the same status).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- 1. Flip the status text everywhere it appears -----------------------
$ov = $wb.Worksheets.Item("Overview")
foreach ($addr in @("B2", "C2", "B3", "C3")) {
    if ($ov.Range($addr).Text -eq $oldStatus) {
        $ov.Range($addr).Value = $newStatus
    }
}

$langSheets = @("zh-cn", "de-de")
foreach ($name in $langSheets) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in @(2, 3)) {
        $c = $ws.Cells.Item($row, 3)
        if ($c.Text -eq $oldStatus) {
            $c.Value = $newStatus
        }
    }
}

# ---- helper: find the Address of the hyperlink anchored at a given cell --
function Get-LinkAddress($ws, $cellAddr) {
    foreach ($link in $ws.Hyperlinks) {
        if ($link.Range.Address() -eq $cellAddr) {
            return $link.Address
        }
    }
    return $null
}

# ---- helper: mirror one existing link cell into a brand-new link cell ----
function Copy-LinkCell($ws, $srcRow, $srcCol, $dstRow, $dstCol) {
    $src = $ws.Cells.Item($srcRow, $srcCol)
    $dst = $ws.Cells.Item($dstRow, $dstCol)
    $srcAddr = $src.Address()
    $target = Get-LinkAddress $ws $srcAddr
    $display = $src.Text

    $ws.Hyperlinks.Add($dst, $target, "", "", $display) | Out-Null

    # Match the look of the existing hyperlink columns (underline + the
    # workbook's custom link color).
    $dst.Font.Underline = 2
    $dst.Font.Color = 15570276
}

# ---- 2. Handback datetimes, one per language sheet ------------------------
$handbackTimes = @{ "zh-cn" = "2016-03-22 06:22:04"; "de-de" = "2016-03-22 06:22:10" }

foreach ($name in $langSheets) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($row in @(2, 3)) {
        # F = Latest Target File  <- mirrors A = Source File Name
        Copy-LinkCell $ws $row 1 $row 6
        # G = Latest Handback File <- mirrors D = Latest Handoff File
        Copy-LinkCell $ws $row 4 $row 7

        # H = Latest Handback DateTime
        $ws.Cells.Item($row, 8).Value = $handbackTimes[$name]
    }
}
